$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine how many data rows currently exist (row 1 is the header).
$lastRow = $ws.UsedRange.Rows.Count

# New header cells: AD1=Wins, AE1=Losses, AF1=Ties
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the formatting of the existing header cell (AC1) onto the new
# header cells so they look consistent with the rest of row 1
# (bold, centered, bordered).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the season record (Wins/Losses/Ties) for every data row.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 74
    $ws.Cells.Item($r, 31).Value = 88
    $ws.Cells.Item($r, 32).Value = 0
}
